$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    3  = @("特变电工", "特变电工", "特变电工")
    4  = @("海马汽车", "吉视传媒", "吉视传媒")
    5  = @("雪人集团", "福龙马",   "雪人集团")
    6  = @("福龙马",   "雪人集团", "海马汽车")
    7  = @("阳光电源", "海马汽车", "福龙马")
    8  = @("合富中国", "神马电力", "合富中国")
    9  = @("海南发展", "阳光电源", "神马电力")
    10 = @("中能电气", "中能电气", "中钨高新")
    11 = @("吉视传媒", "双杰电气", "海南发展")
    12 = @("神马电力", "海南发展", "盈新发展")
    13 = @("多氟多",   "海陆重工", "多氟多")
    14 = @("粤 传 媒", "万里马",   "阳光电源")
    15 = @("海陆重工", "保变电气", "粤传媒")
    16 = @("双杰电气", "多氟多",   "漳州发展")
    17 = @("盈新发展", "粤 传 媒", "安泰集团")
    18 = @("海峡股份", "合富中国", "大连圣亚")
    19 = @("京泉华",   "中国电影", "中能电气")
    20 = @("金盘科技", "工业富联", "海陆重工")
    21 = @("安泰集团", "安泰集团", "神州信息")
}

foreach ($row in $values.Keys) {
    $rowVals = $values[$row]
    $ws.Cells.Item($row, 1).Value = $rowVals[0]
    $ws.Cells.Item($row, 2).Value = $rowVals[1]
    $ws.Cells.Item($row, 3).Value = $rowVals[2]
}
